# Budget.xlsx update:
#  - Logged a new week of labor hours (Labor Hours sheet), including the
#    Feb 2 presentation day, which ripples into the Overall Budget totals.
#  - Corrected the oscilloscope rental cost estimate (Material Costs!B5:
#    3.3 -> 13.3), which ripples into Overall Budget Direct Material Costs.
#  - Updated sheet selections / active tab to reflect where the author
#    was last working (Labor Hours).

$wb = $excel.ActiveWorkbook

# --- Material Costs sheet: correct the oscilloscope rental cost ---
$wsMat = $wb.Worksheets.Item("Material Costs")
$wsMat.Range("B5").Formula = "=13.3"

# --- Labor Hours sheet: fill in the new week of logged hours ---
$wsLabor = $wb.Worksheets.Item("Labor Hours")

# Monday column (G) had been left blank for every earlier week; zero it out.
$wsLabor.Range("G2").Value = 0
$wsLabor.Range("G3").Value = 0
$wsLabor.Range("G4").Value = 0
$wsLabor.Range("G5").Value = 0
$wsLabor.Range("G6").Value = 0
$wsLabor.Range("G7").Value = 0
$wsLabor.Range("G8").Value = 0
$wsLabor.Range("G9").Value = 0
$wsLabor.Range("G10").Value = 0
$wsLabor.Range("G11").Value = 0
$wsLabor.Range("G12").Value = 0
$wsLabor.Range("G13").Value = 0
$wsLabor.Range("G14").Value = 0
$wsLabor.Range("G15").Value = 0
$wsLabor.Range("G16").Value = 0

# Week of Jan 26 - Feb 1 (rows 17-23): Dirk/Erik/Mohammed hours logged,
# including the Feb 2 Presentation prep days.
$wsLabor.Range("B17").Value = 0
$wsLabor.Range("C17").Value = 0
$wsLabor.Range("D17").Value = 0
$wsLabor.Range("E17").Value = 0
$wsLabor.Range("F17").Value = 0
$wsLabor.Range("G17").Value = 0

$wsLabor.Range("B18").Value = 0
$wsLabor.Range("C18").Value = 0
$wsLabor.Range("D18").Value = 0
$wsLabor.Range("E18").Value = 0
$wsLabor.Range("F18").Value = 0
$wsLabor.Range("G18").Value = 0

$wsLabor.Range("B19").Value = 3
$wsLabor.Range("C19").Value = 3
$wsLabor.Range("D19").Value = 3
$wsLabor.Range("E19").Value = 0
$wsLabor.Range("F19").Value = 0
$wsLabor.Range("G19").Value = 0

$wsLabor.Range("B20").Value = 0
$wsLabor.Range("C20").Value = 0
$wsLabor.Range("D20").Value = 0
$wsLabor.Range("E20").Value = 0
$wsLabor.Range("F20").Value = 0
$wsLabor.Range("G20").Value = 0

$wsLabor.Range("B21").Value = 3
$wsLabor.Range("C21").Value = 3
$wsLabor.Range("D21").Value = 3
$wsLabor.Range("E21").Value = 0
$wsLabor.Range("F21").Value = 0
$wsLabor.Range("G21").Value = 0

$wsLabor.Range("B22").Value = 0
$wsLabor.Range("C22").Value = 0
$wsLabor.Range("D22").Value = 0
$wsLabor.Range("E22").Value = 0
$wsLabor.Range("F22").Value = 0
$wsLabor.Range("G22").Value = 0

$wsLabor.Range("B23").Value = 6
$wsLabor.Range("C23").Value = 6
$wsLabor.Range("D23").Value = 3
$wsLabor.Range("E23").Value = 0
$wsLabor.Range("F23").Value = 0
$wsLabor.Range("G23").Value = 0

# --- Selections / active sheet bookkeeping, matching the author's last view ---
$wsOverall = $wb.Worksheets.Item("Overall Budget")
[void]$wsOverall.Activate()
$wsOverall.Range("M13").Select() | Out-Null

[void]$wsMat.Activate()
$wsMat.Range("B6").Select() | Out-Null

[void]$wsLabor.Activate()
$wsLabor.Range("G24").Select() | Out-Null
